$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06440233333333333
$ws.Range("H2").Value = 0.193207
$ws.Range("I2").Value = 0.03647206354366116
$ws.Range("J2").Value = 0.03647206354366116
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 0.9223625541223331
$ws.Range("R2").Value = 8.301262987101
$ws.Range("S2").Value = 0.01075768744420563
$ws.Range("T2").Value = 0.01075768744420563

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06440233333333333
$ws.Range("H3").Value = 0.193207
$ws.Range("I3").Value = 0.03647206354366116
$ws.Range("J3").Value = 0.03647206354366116
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.084169
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 1.744283679994333
$ws.Range("R3").Value = 15.698553119949
$ws.Range("S3").Value = 0.02034390767442094
$ws.Range("T3").Value = 0.02034390767442094

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06440233333333333
$ws.Range("H4").Value = 0.193207
$ws.Range("I4").Value = 0.03647206354366116
$ws.Range("J4").Value = 0.03647206354366116
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 0.4604631803107777
$ws.Range("R4").Value = 4.144168622796999
$ws.Range("S4").Value = 0.005370468425034588
$ws.Range("T4").Value = 0.005370468425034587

# Row 5
$ws.Range("I5").Value = 0.8194013021867156
$ws.Range("J5").Value = 0.8194013021867155
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 20.722300975138
$ws.Range("R5").Value = 186.500708776242
$ws.Range("S5").Value = 0.2416880824345843
$ws.Range("T5").Value = 0.2416880824345842

# Row 6
$ws.Range("I6").Value = 0.8194013021867156
$ws.Range("J6").Value = 0.8194013021867155
$ws.Range("M6").Value = 27.084169
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("Q6").Value = 39.18802995776201
$ws.Range("R6").Value = 352.6922696198581
$ws.Range("S6").Value = 0.4570573425337225
$ws.Range("T6").Value = 0.4570573425337225

# Row 7
$ws.Range("I7").Value = 0.8194013021867156
$ws.Range("J7").Value = 0.8194013021867155
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 10.34501733371933
$ws.Range("R7").Value = 93.10515600347399
$ws.Range("S7").Value = 0.1206558772184087
$ws.Range("T7").Value = 0.1206558772184087

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2544986666666667
$ws.Range("H8").Value = 0.763496
$ws.Range("I8").Value = 0.1441266342696234
$ws.Range("J8").Value = 0.1441266342696234
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 3.644899618658666
$ws.Range("R8").Value = 32.804096567928
$ws.Range("S8").Value = 0.0425111477995167
$ws.Range("T8").Value = 0.0425111477995167

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2544986666666667
$ws.Range("H9").Value = 0.763496
$ws.Range("I9").Value = 0.1441266342696234
$ws.Range("J9").Value = 0.1441266342696234
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.084169
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 6.892884898274667
$ws.Range("R9").Value = 62.035964084472
$ws.Range("S9").Value = 0.08039300922735557
$ws.Range("T9").Value = 0.08039300922735557

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.2544986666666667
$ws.Range("H10").Value = 0.763496
$ws.Range("I10").Value = 0.1441266342696234
$ws.Range("J10").Value = 0.1441266342696234
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 1.819612106779555
$ws.Range("R10").Value = 16.376508961016
$ws.Range("S10").Value = 0.02122247724275108
$ws.Range("T10").Value = 0.02122247724275108

